# Update "total_risk" (column R) and "total_risk_resp" (column S) figures
# on the facility-demographics sheet with the latest AirToxics NATA data.
# (output updates for allocation and transitions rule with newest airtoxics nata data)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Linde-Decatur
$ws.Range("R2").Value = 42.0512820512821
$ws.Range("S2").Value = 0.471794871794872

# Row 3 - A-GAS
$ws.Range("R3").Value = 46.6666666666667
$ws.Range("S3").Value = 0.5

# Row 4 - AEROPRES-SANDIMAS
$ws.Range("R4").Value = 30.2912621359223
$ws.Range("S4").Value = 0.399514563106796

# Row 5 - CALAMCO
$ws.Range("R5").Value = 29.8780487804878
$ws.Range("S5").Value = 0.347560975609756

# Row 6 - Linde-Whiting
$ws.Range("R6").Value = 20.6060606060606
$ws.Range("S6").Value = 0.301212121212121

# Row 7 - Diversified-CPC
$ws.Range("R7").Value = 20
$ws.Range("S7").Value = 0.27

# Row 9 - APC-Geismar
$ws.Range("R9").Value = 61.1764705882353
$ws.Range("S9").Value = 0.405882352941176

# Row 10 - Honeywell-Geismar
$ws.Range("R10").Value = 61.4285714285714
$ws.Range("S10").Value = 0.45

# Row 11 - APC-PortAuthur
$ws.Range("R11").Value = 50
$ws.Range("S11").Value = 0.3

# Row 12 - AEROPRES-SIBLEY
$ws.Range("R12").Value = 33.8095238095238
$ws.Range("S12").Value = 0.414285714285714

# Row 13 - HaltermanCarless (total_risk unchanged, only total_risk_resp changes)
$ws.Range("S13").Value = 0.314285714285714

# Row 14 - Chemours-CorpusChristie
$ws.Range("R14").Value = 19.047619047619
$ws.Range("S14").Value = 0.19047619047619

# Row 15 - DiversifiedG&O (total_risk unchanged, only total_risk_resp changes)
$ws.Range("S15").Value = 0.31
